# Standardise pluralisation of names across the workbook.
$wb = $excel.ActiveWorkbook

# --- Colors sheet and pop_names sheet: pluralise the Type.Name labels ---
# (new text is introduced in this order so the shared-string table is
# rebuilt in the same order as the canonical edit)
$colors   = $wb.Worksheets.Item("Colors")
$popNames = $wb.Worksheets.Item("pop_names")

$colors.Range("A10").Value = "American Wolves"
$popNames.Range("C74").Value = "American Wolves"

$colors.Range("A11").Value = "Eurasian Wolves"
$popNames.Range("C75").Value = "Eurasian Wolves"
$popNames.Range("C76").Value = "Eurasian Wolves"
$popNames.Range("C77").Value = "Eurasian Wolves"

$colors.Range("A12").Value = "Ancient Wolves"
$popNames.Range("C78").Value = "Ancient Wolves"

$colors.Range("A13").Value = "Coyotes"
$popNames.Range("C79").Value = "Coyotes"

$colors.Range("A7").Value  = "Dingos"
$popNames.Range("C71").Value = "Dingos"

# --- Restore the selected cell on each sheet as recorded in the edit ---
$popNames.Range("G84").Select()
$colors.Range("A8").Select()
